$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 mirrors row 9's formatting; fill in the new timesheet entry.
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A10").Value = "29.3.2020"
$ws.Range("B10").Value = 0.83333333333333337
$ws.Range("C10").Value = 0.89583333333333337
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("E10").Value = "Pattern Gen 2"
$ws.Range("F10").Value = "RTL and TB"

$ws.Range("A11").Select()
